$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (ID only changes)
$ws.Range("A2").Value = "a0250d49-275b-4792-9531-82a1aa23e91e"

# Update row 3 (ID and email change)
$ws.Range("A3").Value = "898e3b3a-3115-45bd-bbcb-d5bcabeb41c0"
$ws.Range("C3").Value = "joaquim@example.com"

# Update row 4 (ID, name, email all change)
$ws.Range("A4").Value = "9e7b33d3-42c9-4af7-8c12-40da112137f6"
$ws.Range("B4").Value = "Perez"
$ws.Range("C4").Value = "perez@gmail.com"

# Update row 5 (ID, name, email all change)
$ws.Range("A5").Value = "74340479-d55a-41ed-b3d0-89249e2bcc5b"
$ws.Range("B5").Value = "Peraldo"
$ws.Range("C5").Value = "perez2@gmail.com"

# Remove row 6 entirely (shrinks dimension from A1:C6 to A1:C5)
$ws.Rows("6").Delete()
